$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 331.2
$ws.Range("I9").Value = 332.57144
$ws.Range("J9").Value = 328
$ws.Range("K9").Value = 332.57144
$ws.Range("L9").Value = 328
$ws.Range("M9").Value = -163.57144
$ws.Range("N9").Value = -666
$ws.Range("H17").Value = 1307.138
$ws.Range("J17").Value = 1307.138
$ws.Range("L17").Value = 3921.414
$ws.Range("N17").Value = -4257.414
$ws.Range("H19").Value = 986.1111
$ws.Range("I19").Value = 744
$ws.Range("J19").Value = 1055.2858
$ws.Range("K19").Value = 744
$ws.Range("L19").Value = 1055.2858
$ws.Range("M19").Value = -569
$ws.Range("N19").Value = -1405.2858
$ws.Range("H40").Value = 45505770
$ws.Range("I40").Value = 50750
$ws.Range("J40").Value = 55606880
$ws.Range("K40").Value = 50750
$ws.Range("L40").Value = 55606880
$ws.Range("M40").Value = -50575
$ws.Range("N40").Value = -55607230
$ws.Range("H49").Value = 4562.5
$ws.Range("I49").Value = 7000
$ws.Range("J49").Value = 3750
$ws.Range("K49").Value = 21000
$ws.Range("L49").Value = 11250
$ws.Range("M49").Value = -20864
$ws.Range("N49").Value = -11522
$ws.Range("H59").Value = 4737.6
$ws.Range("I59").Value = 2555
$ws.Range("K59").Value = 7665
$ws.Range("M59").Value = -7108
$ws.Range("H80").Value = 747.9375
$ws.Range("I80").Value = 516.36365
$ws.Range("K80").Value = 1549.09095
$ws.Range("M80").Value = -551.09095
$ws.Range("H83").Value = 747.9375
$ws.Range("I83").Value = 516.36365
$ws.Range("K83").Value = 4647.27285
$ws.Range("M83").Value = 344.7271499999997
$ws.Range("H86").Value = 1541990.5
$ws.Range("J86").Value = 3222.2
$ws.Range("L86").Value = 3222.2
$ws.Range("N86").Value = -5468.2
$ws.Range("H89").Value = 1541990.5
$ws.Range("J89").Value = 3222.2
$ws.Range("L89").Value = 16111
$ws.Range("N89").Value = -27343
$ws.Range("H92").Value = 637.7646999999999
$ws.Range("I92").Value = 558.5625
$ws.Range("K92").Value = 558.5625
$ws.Range("M92").Value = 689.4375
$ws.Range("H99").Value = 5435.25
$ws.Range("I99").Value = 1815.6666
$ws.Range("J99").Value = 7607
$ws.Range("K99").Value = 5446.9998
$ws.Range("L99").Value = 22821
$ws.Range("M99").Value = -3948.9998
$ws.Range("N99").Value = -25817
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = ""
$ws.Range("H133").Value = 93780
$ws.Range("J133").Value = 93780
$ws.Range("L133").Value = 93780
$ws.Range("N133").Value = -103900
$ws.Range("H137").Value = 3167.8635
$ws.Range("J137").Value = 6707.9414
$ws.Range("L137").Value = 20123.8242
$ws.Range("N137").Value = -25223.8242
$ws.Range("H138").Value = 3433.4407
$ws.Range("I138").Value = 1588.8572
$ws.Range("J138").Value = 4007.311
$ws.Range("K138").Value = 4766.571599999999
$ws.Range("L138").Value = 12021.933
$ws.Range("M138").Value = 373.4284000000007
$ws.Range("N138").Value = -22301.933
$ws.Range("H141").Value = 6810.9585
$ws.Range("I141").Value = 6025.1055
$ws.Range("K141").Value = 18075.3165
$ws.Range("M141").Value = -12895.3165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3506.5715
$ws.Range("I32").Value = 1802.3077
$ws.Range("K32").Value = 1802.3077
$ws.Range("M32").Value = -1515.3077
$ws.Range("H61").Value = 36102
$ws.Range("I61").Value = 36102
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 36102
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -35890
$ws.Range("N61").Value = ""
$ws.Range("H74").Value = 7354398.5
$ws.Range("I74").Value = 11905727
$ws.Range("K74").Value = 11905727
$ws.Range("M74").Value = -11904853
$ws.Range("H77").Value = 7354398.5
$ws.Range("I77").Value = 11905727
$ws.Range("K77").Value = 59528635
$ws.Range("M77").Value = -59524267
$ws.Range("H102").Value = 335309.38
$ws.Range("I102").Value = 623601.4
$ws.Range("K102").Value = 623601.4
$ws.Range("M102").Value = -621979.4
$ws.Range("H136").Value = 36102
$ws.Range("I136").Value = 36102
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 108306
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -105756
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1319.1875
$ws.Range("I107").Value = 1229.5
$ws.Range("J107").Value = 1707.8334
$ws.Range("K107").Value = 1229.5
$ws.Range("L107").Value = 1707.8334
$ws.Range("M107").Value = 690.5
$ws.Range("N107").Value = -5547.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 19523.5
$ws.Range("I36").Value = 19047
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 19047
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -18659
$ws.Range("N36").Value = -20776
$ws.Range("H40").Value = 19523.5
$ws.Range("I40").Value = 19047
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 19047
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = -18887
$ws.Range("N40").Value = -20320
$ws.Range("H62").Value = 30090.908
$ws.Range("I62").Value = 2866.3333
$ws.Range("J62").Value = 40300.125
$ws.Range("K62").Value = 2866.3333
$ws.Range("L62").Value = 40300.125
$ws.Range("M62").Value = -2242.3333
$ws.Range("N62").Value = -41548.125
$ws.Range("H65").Value = 30090.908
$ws.Range("I65").Value = 2866.3333
$ws.Range("J65").Value = 40300.125
$ws.Range("K65").Value = 14331.6665
$ws.Range("L65").Value = 201500.625
$ws.Range("M65").Value = -11211.6665
$ws.Range("N65").Value = -207740.625
$ws.Range("H86").Value = 3499.8
$ws.Range("I86").Value = 3833
$ws.Range("K86").Value = 3833
$ws.Range("M86").Value = -2710
$ws.Range("H89").Value = 3499.8
$ws.Range("I89").Value = 3833
$ws.Range("K89").Value = 19165
$ws.Range("M89").Value = -13549
$ws.Range("H99").Value = 13046.733
$ws.Range("I99").Value = 20530.143
$ws.Range("J99").Value = 6498.75
$ws.Range("K99").Value = 20530.143
$ws.Range("L99").Value = 6498.75
$ws.Range("M99").Value = -19032.143
$ws.Range("N99").Value = -9494.75
$ws.Range("H107").Value = 1136749.8
$ws.Range("I107").Value = 1515395.9
$ws.Range("K107").Value = 1515395.9
$ws.Range("M107").Value = -1513475.9
$ws.Range("H126").Value = 13046.733
$ws.Range("I126").Value = 20530.143
$ws.Range("J126").Value = 6498.75
$ws.Range("K126").Value = 61590.429
$ws.Range("L126").Value = 19496.25
$ws.Range("M126").Value = -59120.429
$ws.Range("N126").Value = -24436.25
$ws.Range("H132").Value = 10103169
$ws.Range("I132").Value = 13335339
$ws.Range("K132").Value = 40006017
$ws.Range("M132").Value = -40003487
$ws.Range("H141").Value = 92122
$ws.Range("J141").Value = 100893.336
$ws.Range("L141").Value = 100893.336
$ws.Range("N141").Value = -111253.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 389.8
$ws.Range("J33").Value = 549.5
$ws.Range("L33").Value = 3297
$ws.Range("N33").Value = -3863
$ws.Range("H40").Value = 1427.25
$ws.Range("I40").Value = 97.875
$ws.Range("J40").Value = 2756.625
$ws.Range("K40").Value = 391.5
$ws.Range("L40").Value = 11026.5
$ws.Range("M40").Value = -322.5
$ws.Range("N40").Value = -11164.5
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -35242
$ws.Range("H107").Value = 1004.2963
$ws.Range("I107").Value = 1417.25
$ws.Range("J107").Value = 830.4211
$ws.Range("K107").Value = 4251.75
$ws.Range("L107").Value = 2491.2633
$ws.Range("M107").Value = -2331.75
$ws.Range("N107").Value = -6331.263300000001
$ws.Range("H113").Value = 1073.3334
$ws.Range("J113").Value = 1211
$ws.Range("L113").Value = 3633
$ws.Range("N113").Value = -7973

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 10000
$ws.Range("K30").Value = 10000
$ws.Range("M30").Value = -9892
$ws.Range("H38").Value = 18750
$ws.Range("J38").Value = 18750
$ws.Range("L38").Value = 18750
$ws.Range("N38").Value = -19570
$ws.Range("H61").Value = 3798
$ws.Range("I61").Value = 2600
$ws.Range("K61").Value = 2600
$ws.Range("M61").Value = -2398
$ws.Range("H113").Value = 3798
$ws.Range("I113").Value = 2600
$ws.Range("K113").Value = 2600
$ws.Range("M113").Value = -430
$ws.Range("H122").Value = 47622210
$ws.Range("J122").Value = 28576226
$ws.Range("L122").Value = 85728678
$ws.Range("N122").Value = -85733578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14623
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = ""
$ws.Range("H81").Value = 3795714.8
$ws.Range("I81").Value = 2612541.5
$ws.Range("K81").Value = 5225083
$ws.Range("M81").Value = -5224022
$ws.Range("H84").Value = 3795714.8
$ws.Range("I84").Value = 2612541.5
$ws.Range("K84").Value = 26125415
$ws.Range("M84").Value = -26120111
$ws.Range("H100").Value = 1178418
$ws.Range("I100").Value = 2001323.8
$ws.Range("K100").Value = 4002647.6
$ws.Range("M100").Value = -4002106.6
$ws.Range("H136").Value = 8021.6313
$ws.Range("I136").Value = 3969.3103
$ws.Range("K136").Value = 11907.9309
$ws.Range("M136").Value = -9357.930899999999
